# BusinessTOC.docx edit script
#
# 1) In the "Rule : Statement (...)" metamodel-API line, swap the Class/Kind
#    roles of the "subject" and "condition" parameters.
# 2) After the "Functional transform:" paragraph, add a new blank paragraph
#    followed by a new paragraph of explanatory text (mirroring the layout
#    already used elsewhere in the doc, e.g. under "Functional apply:" /
#    "Functional query:").

$d = $word.ActiveDocument

# --- Edit 1: Class/Kind swap in the Rule Statement signature -------------
$d.Content.Find.Execute(
    "Rule : Statement (Rule, subject : Class, condition : Kind, transform : Flow);",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Rule : Statement (Rule, subject : Kind, condition : Class, transform : Flow);",
    2)

# --- Edit 2: insert two new paragraphs right after "Functional transform:" -
$rng = $d.Content
$rng.Find.Execute("Functional transform:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$anchorPara = $rng.Paragraphs(1)

# New blank paragraph directly after "Functional transform:".
$anchorPara.Range.InsertParagraphAfter()
$blankPara = $anchorPara.Next()

# New paragraph with the explanatory sentence, directly after the blank one.
$blankPara.Range.InsertParagraphAfter()
$textPara = $blankPara.Next()
$textPara.Range.InsertBefore("In the context of a Facts dialog given matching Concepts a Template could be matched which activates a Rule Flow (pattern / transform) which updates players LHS with RHS.")
